# Auto-generated Excel COM-interop script to apply scheduled-runner updates
# to cached market price / profit values across the Sheets workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2096.875
$ws.Range("I40").Value = 2065.3845
$ws.Range("J40").Value = 2233.3333
$ws.Range("K40").Value = 2065.3845
$ws.Range("L40").Value = 2233.3333
$ws.Range("M40").Value = -1890.3845
$ws.Range("N40").Value = -2583.3333

$ws.Range("H64").Value = 3576
$ws.Range("I64").Value = 3361.2122
$ws.Range("J64").Value = 4166.6665
$ws.Range("K64").Value = 3361.2122
$ws.Range("L64").Value = 4166.6665
$ws.Range("M64").Value = -3113.2122
$ws.Range("N64").Value = -4662.6665

$ws.Range("H67").Value = 3576
$ws.Range("I67").Value = 3361.2122
$ws.Range("J67").Value = 4166.6665
$ws.Range("K67").Value = 3361.2122
$ws.Range("L67").Value = 4166.6665
$ws.Range("M67").Value = -2503.2122
$ws.Range("N67").Value = -5882.6665

$ws.Range("H74").Value = 6281.353
$ws.Range("I74").Value = 3630.3333
$ws.Range("J74").Value = 7727.364
$ws.Range("K74").Value = 3630.3333
$ws.Range("L74").Value = 7727.364
$ws.Range("M74").Value = -2694.3333
$ws.Range("N74").Value = -9599.364

$ws.Range("H76").Value = 3432.195
$ws.Range("I76").Value = 2990.9092
$ws.Range("K76").Value = 2990.9092
$ws.Range("M76").Value = -2675.9092

$ws.Range("H77").Value = 6281.353
$ws.Range("I77").Value = 3630.3333
$ws.Range("J77").Value = 7727.364
$ws.Range("K77").Value = 18151.6665
$ws.Range("L77").Value = 38636.82
$ws.Range("M77").Value = -13471.6665
$ws.Range("N77").Value = -47996.82

$ws.Range("H79").Value = 3432.195
$ws.Range("I79").Value = 2990.9092
$ws.Range("K79").Value = 2990.9092
$ws.Range("M79").Value = -1898.9092

$ws.Range("H87").Value = 15364.389
$ws.Range("J87").Value = 15364.389
$ws.Range("L87").Value = 15364.389
$ws.Range("N87").Value = -17860.389

$ws.Range("H90").Value = 15364.389
$ws.Range("J90").Value = 15364.389
$ws.Range("L90").Value = 46093.167
$ws.Range("N90").Value = -58573.167

$ws.Range("H129").Value = 1020.9359
$ws.Range("J129").Value = 1040.9865
$ws.Range("L129").Value = 3122.9595
$ws.Range("N129").Value = -13122.9595

$ws.Range("H137").Value = 1306.6154
$ws.Range("I137").Value = 1579.4546
$ws.Range("J137").Value = 1106.5333
$ws.Range("K137").Value = 4738.3638
$ws.Range("L137").Value = 3319.5999
$ws.Range("M137").Value = -2188.3638
$ws.Range("N137").Value = -8419.599900000001

$ws.Range("H138").Value = 4516.173
$ws.Range("I138").Value = 2299.862
$ws.Range("J138").Value = 5752.1924
$ws.Range("K138").Value = 6899.586
$ws.Range("L138").Value = 17256.5772
$ws.Range("M138").Value = -1759.586
$ws.Range("N138").Value = -27536.5772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2016.5483
$ws.Range("I61").Value = 1760.1364
$ws.Range("J61").Value = 2643.3333
$ws.Range("K61").Value = 1760.1364
$ws.Range("L61").Value = 2643.3333
$ws.Range("M61").Value = -1548.1364
$ws.Range("N61").Value = -3067.3333

$ws.Range("H63").Value = 2572
$ws.Range("I63").Value = 2000.8
$ws.Range("K63").Value = 2000.8
$ws.Range("M63").Value = -1314.8

$ws.Range("H66").Value = 2572
$ws.Range("I66").Value = 2000.8
$ws.Range("K66").Value = 10004
$ws.Range("M66").Value = -6572

$ws.Range("H132").Value = 2065.5232
$ws.Range("I132").Value = 1653.3096
$ws.Range("K132").Value = 4959.9288
$ws.Range("M132").Value = -2429.9288

$ws.Range("H136").Value = 2016.5483
$ws.Range("I136").Value = 1760.1364
$ws.Range("J136").Value = 2643.3333
$ws.Range("K136").Value = 5280.4092
$ws.Range("L136").Value = 7929.999899999999
$ws.Range("M136").Value = -2730.4092
$ws.Range("N136").Value = -13029.9999

$ws.Range("H139").Value = 67202.78
$ws.Range("J139").Value = 67202.78
$ws.Range("L139").Value = 67202.78
$ws.Range("N139").Value = -77482.78

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2119.087
$ws.Range("I105").Value = 1998.1177
$ws.Range("K105").Value = 1998.1177
$ws.Range("M105").Value = -251.1177

$ws.Range("H134").Value = 1594.1754
$ws.Range("I134").Value = 1209.6923
$ws.Range("J134").Value = 2427.2222
$ws.Range("K134").Value = 3629.0769
$ws.Range("L134").Value = 7281.6666
$ws.Range("M134").Value = -1094.0769
$ws.Range("N134").Value = -12351.6666

$ws.Range("H140").Value = 58031.668
$ws.Range("J140").Value = 58031.668
$ws.Range("L140").Value = 58031.668
$ws.Range("N140").Value = -68391.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4336.8076
$ws.Range("I31").Value = 2139.2856
$ws.Range("J31").Value = 9930.5
$ws.Range("K31").Value = 2139.2856
$ws.Range("L31").Value = 9930.5
$ws.Range("M31").Value = -1844.2856
$ws.Range("N31").Value = -10520.5

$ws.Range("H34").Value = 4336.8076
$ws.Range("I34").Value = 2139.2856
$ws.Range("J34").Value = 9930.5
$ws.Range("K34").Value = 2139.2856
$ws.Range("L34").Value = 9930.5
$ws.Range("M34").Value = -1937.2856
$ws.Range("N34").Value = -10334.5

$ws.Range("H62").Value = 671074.4
$ws.Range("I62").Value = 4502.5
$ws.Range("J62").Value = 913464.2
$ws.Range("K62").Value = 4502.5
$ws.Range("L62").Value = 913464.2
$ws.Range("M62").Value = -3878.5
$ws.Range("N62").Value = -914712.2

$ws.Range("H65").Value = 671074.4
$ws.Range("I65").Value = 4502.5
$ws.Range("J65").Value = 913464.2
$ws.Range("K65").Value = 22512.5
$ws.Range("L65").Value = 4567321
$ws.Range("M65").Value = -19392.5
$ws.Range("N65").Value = -4573561

$ws.Range("H108").Value = 22054.5
$ws.Range("I108").Value = 9999
$ws.Range("J108").Value = 34110
$ws.Range("K108").Value = 9999
$ws.Range("L108").Value = 34110
$ws.Range("M108").Value = -6159
$ws.Range("N108").Value = -41790

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 63
$ws.Range("J12").Value = 60.266666
$ws.Range("L12").Value = 180.799998
$ws.Range("N12").Value = -526.799998

$ws.Range("H33").Value = 90.666664
$ws.Range("I33").Value = 86
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 516
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = -233
$ws.Range("N33").Value = -1166

$ws.Range("H131").Value = 842.9697
$ws.Range("I131").Value = 406.15384
$ws.Range("J131").Value = 909
$ws.Range("K131").Value = 1218.46152
$ws.Range("L131").Value = 2727
$ws.Range("M131").Value = 3821.53848
$ws.Range("N131").Value = -12807

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6379.9
$ws.Range("I70").Value = 4983.636
$ws.Range("J70").Value = 8086.4443
$ws.Range("K70").Value = 4983.636
$ws.Range("L70").Value = 8086.4443
$ws.Range("M70").Value = -4713.636
$ws.Range("N70").Value = -8626.444299999999

$ws.Range("H73").Value = 6379.9
$ws.Range("I73").Value = 4983.636
$ws.Range("J73").Value = 8086.4443
$ws.Range("K73").Value = 4983.636
$ws.Range("L73").Value = 8086.4443
$ws.Range("M73").Value = -4047.636
$ws.Range("N73").Value = -9958.444299999999

$ws.Range("H80").Value = 2416.8333
$ws.Range("I80").Value = 2447.5
$ws.Range("J80").Value = 2401.5
$ws.Range("K80").Value = 2447.5
$ws.Range("L80").Value = 2401.5
$ws.Range("M80").Value = -1449.5
$ws.Range("N80").Value = -4397.5

$ws.Range("H83").Value = 2416.8333
$ws.Range("I83").Value = 2447.5
$ws.Range("J83").Value = 2401.5
$ws.Range("K83").Value = 12237.5
$ws.Range("L83").Value = 12007.5
$ws.Range("M83").Value = -7245.5
$ws.Range("N83").Value = -21991.5

$ws.Range("H122").Value = 58741.16
$ws.Range("I122").Value = 72815.86
$ws.Range("K122").Value = 218447.58
$ws.Range("M122").Value = -215997.58

$ws.Range("H135").Value = 41407.777
$ws.Range("J135").Value = 41407.777
$ws.Range("L135").Value = 41407.777
$ws.Range("N135").Value = -51547.777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3360.8
$ws.Range("I100").Value = 2333.3333
$ws.Range("J100").Value = 4902
$ws.Range("K100").Value = 2333.3333
$ws.Range("L100").Value = 4902
$ws.Range("M100").Value = -1792.3333
$ws.Range("N100").Value = -5984

$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774

$ws.Range("H127").Value = 57182
$ws.Range("J127").Value = 57182
$ws.Range("L127").Value = 57182
$ws.Range("N127").Value = -67102

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H132").Value = 2869.457
$ws.Range("I132").Value = 2595.5
$ws.Range("J132").Value = 3159.5293
$ws.Range("K132").Value = 7786.5
$ws.Range("L132").Value = 9478.5879
$ws.Range("M132").Value = -5256.5
$ws.Range("N132").Value = -14538.5879

$ws.Range("H137").Value = 56773.625
$ws.Range("J137").Value = 56773.625
$ws.Range("L137").Value = 56773.625
$ws.Range("N137").Value = -66973.625
